$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything out - the sheet content is being fully replaced.
$ws.Cells.Clear()

# Row 1 - header row (bold)
$ws.Range("A1").Value = "Column"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Required"
$ws.Range("D1").Value = "Source"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Example"
$ws.Range("A1:F1").Font.Bold = $true

# Row 2 - management_unit
$ws.Range("A2").Value = "management_unit"
$ws.Range("B2").Value = "As is"
$ws.Range("C2").Value = "In"
$ws.Range("D2").Value = "Project Data"
$ws.Range("E2").Value = "This column describes the name of the management unit used for the geographic region (NLP 2018), or the organisation unit responsible for activities in the area (but can reach into other geographic regions)  "
$ws.Range("F2").Value = "Condamine"

# Row 3 - grant_id
$ws.Range("A3").Value = "grant_id"
$ws.Range("B3").Value = "As is"
$ws.Range("C3").Value = "In"
$ws.Range("D3").Value = "Project Data"
$ws.Range("E3").Value = "This column describes the human readable unique ID assigned to a project"
$ws.Range("F3").Value = "RLP-MU46-P2"

# Row 4 - activity_id
$ws.Range("A4").Value = "activity_id"
$ws.Range("B4").Value = "As is"
$ws.Range("C4").Value = "Out"
$ws.Range("D4").Value = "Project Data"
$ws.Range("E4").Value = "N/A"

# Row 5 - project_id
$ws.Range("A5").Value = "project_id"
$ws.Range("B5").Value = "As is"
$ws.Range("C5").Value = "Out"
$ws.Range("D5").Value = "Project Data"
$ws.Range("E5").Value = "N/A"

# Row 6 - program
$ws.Range("A6").Value = "program"
$ws.Range("B6").Value = "As is"
$ws.Range("C6").Value = "In"
$ws.Range("D6").Value = "Project Data"
$ws.Range("E6").Value = "This column describes the program under which the project is being conducted (i.e. source of funding)"
$ws.Range("F6").Value = "National Landcare Program"

# Row 7 - sub_program
$ws.Range("A7").Value = "sub_program"
$ws.Range("B7").Value = "As is"
$ws.Range("C7").Value = "In"
$ws.Range("D7").Value = "Project Data"
$ws.Range("E7").Value = "This column describes the sub-program under which the project is being conducted "
$ws.Range("F7").Value = "Regional Land Partnerships"

# New column F needs a best-fit-style custom width, matching the pattern of
# the pre-existing bestFit columns A/B.
$ws.Columns.Item(6).ColumnWidth = 23

# Move the active selection, matching the saved cursor position.
$null = $ws.Range("E12").Select()
